$d = $word.ActiveDocument

# 1. PERSON_34 paragraph: "s Ritą Holasovou" -> "s Ritą [[PERSON_34]]"
$d.Content.Find.Execute(
    "[[PERSON_34]] – „o [[PERSON_34]]“, „s Ritą Holasovou“",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[[PERSON_34]] – „o [[PERSON_34]]“, „s Ritą [[PERSON_34]]“", 2)

# 2. PERSON_56 paragraph: fix "k [[PERSON_57]]" -> "k [[PERSON_56]]" and add a
#    brand-new PERSON_57 paragraph right after it.
$d.Content.Find.Execute(
    "[[PERSON_56]] – „bez [[PERSON_56]]“, „k [[PERSON_57]]“",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[[PERSON_56]] – „bez [[PERSON_56]]“, „k [[PERSON_56]]“", 2)

$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -match "\[\[PERSON_56\]\]") {
        $para.Range.InsertParagraphAfter()
        $d.Paragraphs($i + 1).Range.Text = "[[PERSON_57]] – „bez [[PERSON_57]]“, „k [[PERSON_57]]“"
        break
    }
}

# 3. PERSON_61 / PERSON_62 paragraphs get relabeled.
$d.Content.Find.Execute(
    "[[PERSON_61]] – „bez [[PERSON_61]]“, „k [[PERSON_61]]“",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[[PERSON_62]] – „bez [[PERSON_62]]“, „k [[PERSON_61]]“", 2)

$d.Content.Find.Execute(
    "[[PERSON_62]] – „bez [[PERSON_62]]“, „k [[PERSON_63]]“",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[[PERSON_63]] – „bez [[PERSON_63]]“, „k [[PERSON_63]]“", 2)

# 4. PERSON_73 / PERSON_74 paragraphs get relabeled.
$d.Content.Find.Execute(
    "[[PERSON_73]] – „bez [[PERSON_73]]“, „k [[PERSON_73]]“",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[[PERSON_75]] – „bez [[PERSON_74]]“, „k [[PERSON_73]]“", 2)

$d.Content.Find.Execute(
    "[[PERSON_74]] – „bez [[PERSON_75]]“, „k [[PERSON_76]]“",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[[PERSON_76]] – „bez [[PERSON_76]]“, „k [[PERSON_76]]“", 2)

# 5. PERSON_84 / PERSON_85 paragraphs get relabeled.
$d.Content.Find.Execute(
    "[[PERSON_84]] – „bez [[PERSON_84]]“, „k [[PERSON_84]]“",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[[PERSON_86]] – „bez [[PERSON_85]]“, „k [[PERSON_84]]“", 2)

$d.Content.Find.Execute(
    "[[PERSON_85]] – „bez [[PERSON_86]]“, „k [[PERSON_87]]“",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[[PERSON_87]] – „bez [[PERSON_87]]“, „k [[PERSON_87]]“", 2)

# 6. PERSON_90 / PERSON_91 paragraphs get relabeled.
$d.Content.Find.Execute(
    "[[PERSON_90]] – „bez [[PERSON_90]]“, „k [[PERSON_90]]“",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[[PERSON_90]] – „bez [[PERSON_91]]“, „k [[PERSON_90]]“", 2)

$d.Content.Find.Execute(
    "[[PERSON_91]] – „bez [[PERSON_92]]“, „k [[PERSON_91]]“",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[[PERSON_92]] – „bez [[PERSON_92]]“, „k [[PERSON_92]]“", 2)

# 7. Remove the trailing PERSON_107 paragraph entirely.
$n = $d.Paragraphs.Count
for ($i = $n; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -match "PERSON_107") {
        $para.Range.Delete()
        break
    }
}

Write-Output "done"
